# Adds a small circled-letter "step" glyph (rendered with the Wingdings 2
# font) in front of four labels on slide 1, and drops the redundant
# "(n) " numeric prefix that used to do the same job as plain text.
#
#   (1) Login request                        -> [u] Login request
#   (2) Auth request                         -> [v] Auth request
#   (3) (11) NGINX Plus exchanges ...         -> [w] (11) NGINX Plus exchanges ...
#   (4) Redirect to original URI              -> [x] Redirect to original URI
#
# Each shape also grows a bit (spAutoFit) to accommodate the extra glyph.

function Get-ShapeById {
    param($slide, [int]$id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

function Add-StepGlyph {
    # Rebuilds the shape's single paragraph as three runs:
    #   1) the Wingdings 2 glyph (e.g. "u" == circled "1")
    #   2) a plain space
    #   3) the (already-formatted) label text
    # and repositions/resizes the shape's box (spAutoFit growth).
    param(
        $shape,
        [string]$glyph,
        [string]$labelText,
        [double]$leftPt,
        [double]$topPt,
        [double]$widthPt,
        [double]$heightPt
    )

    $tr = $shape.TextFrame.TextRange
    $tr.Text = $glyph + " " + $labelText

    $full = $shape.TextFrame2.TextRange

    $iconRange = $full.Characters(1, 1)
    $iconRange.Font.NameAscii = "Wingdings 2"
    $iconRange.Font.NameFarEast = "Malgun Gothic"
    $iconRange.Font.NameComplexScript = "Times New Roman"
    $iconRange.Font.Size = 10

    $spaceRange = $shape.TextFrame2.TextRange.Characters(2, 1)
    $spaceRange.Font.Size = 8

    $shape.Left = $leftPt
    $shape.Top = $topPt
    $shape.Width = $widthPt
    $shape.Height = $heightPt
}

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- (1) Login request ----------------------------------------------------
$sh22 = Get-ShapeById $s 22
$left22 = 4112682 / $EMU_PER_PT
$top22 = 2592066 / $EMU_PER_PT
$w22 = 703719 / $EMU_PER_PT
$h22 = 153888 / $EMU_PER_PT
Add-StepGlyph $sh22 "u" "Login request" $left22 $top22 $w22 $h22

# ---- (2) Auth request -------------------------------------------------------
$sh23 = Get-ShapeById $s 23
$left23 = 5284941 / $EMU_PER_PT
$top23 = 2924180 / $EMU_PER_PT
$w23 = 681277 / $EMU_PER_PT
$h23 = 153888 / $EMU_PER_PT
Add-StepGlyph $sh23 "v" "Auth request" $left23 $top23 $w23 $h23

# ---- (3) (11) NGINX Plus exchanges authorization code for ID / access token -
$sh26 = Get-ShapeById $s 26
$left26 = 6286271 / $EMU_PER_PT
$top26 = 1381299 / $EMU_PER_PT
$w26 = 872212 / $EMU_PER_PT
$h26 = 421013 / $EMU_PER_PT
Add-StepGlyph $sh26 "w" "(11) NGINX Plus exchanges authorization code for ID / access token" $left26 $top26 $w26 $h26

# ---- (4) Redirect to original URI ------------------------------------------
$sh28 = Get-ShapeById $s 28
$left28 = 4047124 / $EMU_PER_PT
$top28 = 2946333 / $EMU_PER_PT
$w28 = 679991 / $EMU_PER_PT
$h28 = 276999 / $EMU_PER_PT
Add-StepGlyph $sh28 "x" "Redirect to original URI" $left28 $top28 $w28 $h28
